# Add season-record columns (Wins / Losses / Ties) to the player stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the last existing header cell (AC1) so the new
# header cells (AD1:AF1) get the same bold/centered/bordered formatting.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row gets the same team season record.
for ($row = 2; $row -le 54; $row++) {
    $ws.Cells.Item($row, 30).Value = 79
    $ws.Cells.Item($row, 31).Value = 83
    $ws.Cells.Item($row, 32).Value = 0
}
